# Update leve-profit metrics across the Anima_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values reflect a refreshed market-price snapshot pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28 (Leve Item ID 27772)
$ws.Range("H28").Value = 16910.445
$ws.Range("I28").Value = 298
$ws.Range("J28").Value = 50135.332
$ws.Range("K28").Value = 298
$ws.Range("L28").Value = 50135.332
$ws.Range("M28").Value = 187
$ws.Range("N28").Value = -51105.332

# Row 34 (Leve Item ID 2160)
$ws.Range("H34").Value = 3183
$ws.Range("J34").Value = 2049
$ws.Range("L34").Value = 2049
$ws.Range("N34").Value = -2455

# Row 36 (Leve Item ID 2160)
$ws.Range("H36").Value = 3183
$ws.Range("J36").Value = 2049
$ws.Range("L36").Value = 2049
$ws.Range("N36").Value = -3479

# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 1500
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 1500
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = -1325
$ws.Range("N40").Value = -1850

# Row 62 (Leve Item ID 27781)
$ws.Range("H62").Value = 2899.4443
$ws.Range("I62").Value = 2073.75
$ws.Range("J62").Value = 3560
$ws.Range("K62").Value = 2073.75
$ws.Range("L62").Value = 3560
$ws.Range("M62").Value = -1449.75
$ws.Range("N62").Value = -4808

# Row 65 (Leve Item ID 27781)
$ws.Range("H65").Value = 2899.4443
$ws.Range("I65").Value = 2073.75
$ws.Range("J65").Value = 3560
$ws.Range("K65").Value = 10368.75
$ws.Range("L65").Value = 17800
$ws.Range("M65").Value = -7248.75
$ws.Range("N65").Value = -24040

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 1540.4789
$ws.Range("I138").Value = 1997.0435
$ws.Range("J138").Value = 1321.7084
$ws.Range("K138").Value = 5991.1305
$ws.Range("L138").Value = 3965.1252
$ws.Range("M138").Value = -851.1305000000002
$ws.Range("N138").Value = -14245.1252

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 1356.2
$ws.Range("I2").Value = 1395.7778
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 1395.7778
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -1282.7778
$ws.Range("N2").Value = -1226

# Row 5 (Leve Item ID 5091)
$ws.Range("H5").Value = 200
$ws.Range("J5").Value = 200
$ws.Range("L5").Value = 200
$ws.Range("N5").Value = -424

# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 2560.3103
$ws.Range("I45").Value = 1846.3889
$ws.Range("J45").Value = 3728.5454
$ws.Range("K45").Value = 1846.3889
$ws.Range("L45").Value = 3728.5454
$ws.Range("M45").Value = -1469.3889
$ws.Range("N45").Value = -4482.5454

# Row 97 (Leve Item ID 19941)
$ws.Range("H97").Value = 767.0625
$ws.Range("I97").Value = 684.86664
$ws.Range("K97").Value = 684.86664
$ws.Range("M97").Value = -188.86664

# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 2384.1428
$ws.Range("I102").Value = 2384.1428
$ws.Range("K102").Value = 2384.1428
$ws.Range("M102").Value = -762.1428000000001

# Row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 2755.5557
$ws.Range("I110").Value = 2600
$ws.Range("J110").Value = 2950
$ws.Range("K110").Value = 2600
$ws.Range("L110").Value = 2950
$ws.Range("M110").Value = -555
$ws.Range("N110").Value = -7040

# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 1356.2
$ws.Range("I116").Value = 1395.7778
$ws.Range("J116").Value = 1000
$ws.Range("K116").Value = 1395.7778
$ws.Range("L116").Value = 1000
$ws.Range("M116").Value = 898.2221999999999
$ws.Range("N116").Value = -5588

# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 1966.5143
$ws.Range("I122").Value = 1670
$ws.Range("J122").Value = 2534.8333
$ws.Range("K122").Value = 5010
$ws.Range("L122").Value = 7604.499899999999
$ws.Range("M122").Value = -2560
$ws.Range("N122").Value = -12504.4999

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 1356.2
$ws.Range("I3").Value = 1395.7778
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 1395.7778
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -1281.7778
$ws.Range("N3").Value = -1228

# Row 4 (Leve Item ID 5091)
$ws.Range("H4").Value = 200
$ws.Range("J4").Value = 200
$ws.Range("L4").Value = 200
$ws.Range("N4").Value = -430

# Row 22 (Leve Item ID 5092)
$ws.Range("H22").Value = 1152.6666
$ws.Range("I22").Value = 1152.6666
$ws.Range("K22").Value = 1152.6666
$ws.Range("M22").Value = -979.6666

# Row 40 (Leve Item ID 19514)
$ws.Range("H40").Value = 95000
$ws.Range("J40").Value = 95000
$ws.Range("L40").Value = 95000
$ws.Range("N40").Value = -95530

# Row 60 (Leve Item ID 43232)
$ws.Range("H60").Value = 29666.666
$ws.Range("J60").Value = 29666.666
$ws.Range("L60").Value = 29666.666
$ws.Range("N60").Value = -30864.666

# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 581.75
$ws.Range("I94").Value = 622.625
$ws.Range("J94").Value = 500
$ws.Range("K94").Value = 622.625
$ws.Range("L94").Value = 500
$ws.Range("M94").Value = -171.625
$ws.Range("N94").Value = -1402

# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 13891288
$ws.Range("I105").Value = 15627274
$ws.Range("K105").Value = 15627274
$ws.Range("M105").Value = -15625527

$ws = $wb.Worksheets.Item("CRP")
# Row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 787.2857
$ws.Range("I105").Value = 700
$ws.Range("K105").Value = 700
$ws.Range("M105").Value = 1047

$ws = $wb.Worksheets.Item("CUL")
# Row 26 (Leve Item ID 4746)
$ws.Range("H26").Value = 394.39285
$ws.Range("I26").Value = 63.285713
$ws.Range("J26").Value = 504.7619
$ws.Range("K26").Value = 189.857139
$ws.Range("L26").Value = 1514.2857
$ws.Range("M26").Value = 98.14286099999998
$ws.Range("N26").Value = -2090.2857

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 5387.964
$ws.Range("I70").Value = 5423.793
$ws.Range("J70").Value = 5348
$ws.Range("K70").Value = 5423.793
$ws.Range("L70").Value = 5348
$ws.Range("M70").Value = -5153.793
$ws.Range("N70").Value = -5888

# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 5387.964
$ws.Range("I73").Value = 5423.793
$ws.Range("J73").Value = 5348
$ws.Range("K73").Value = 5423.793
$ws.Range("L73").Value = 5348
$ws.Range("M73").Value = -4487.793
$ws.Range("N73").Value = -7220

# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 1375.1666
$ws.Range("I113").Value = 965.25
$ws.Range("J113").Value = 2195
$ws.Range("K113").Value = 965.25
$ws.Range("L113").Value = 2195
$ws.Range("M113").Value = 1204.75
$ws.Range("N113").Value = -6535

# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 4186.026
$ws.Range("I122").Value = 1361
$ws.Range("J122").Value = 5951.6665
$ws.Range("K122").Value = 4083
$ws.Range("L122").Value = 17854.9995
$ws.Range("M122").Value = -1633
$ws.Range("N122").Value = -22754.9995

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (Leve Item ID 5289)
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()

# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 13387.8
$ws.Range("I93").Value = 19313.334
$ws.Range("J93").Value = 4499.5
$ws.Range("K93").Value = 19313.334
$ws.Range("L93").Value = 4499.5
$ws.Range("M93").Value = -18065.334
$ws.Range("N93").Value = -6995.5

# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 2326.6667
$ws.Range("I122").Value = 1990
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 5970
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -3520
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("WVR")
# Row 2 (Leve Item ID 3307)
$ws.Range("H2").Value = 2870343
$ws.Range("I2").Value = 10000000
$ws.Range("K2").Value = 10000000
$ws.Range("M2").Value = -9999888

# Row 4 (Leve Item ID 2996)
$ws.Range("H4").Value = 47001.5
$ws.Range("J4").Value = 47001.5
$ws.Range("L4").Value = 47001.5
$ws.Range("N4").Value = -47227.5

# Row 62 (Leve Item ID 12589)
$ws.Range("H62").Value = 52200
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 52200
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 52200
$ws.Range("N62").Value = -53448
$ws.Range("M62").ClearContents()

# Row 65 (Leve Item ID 12589)
$ws.Range("H65").Value = 52200
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 52200
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 261000
$ws.Range("N65").Value = -267240
$ws.Range("M65").ClearContents()

# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 1683.6923
$ws.Range("I126").Value = 1573.5
$ws.Range("J126").Value = 1860
$ws.Range("K126").Value = 4720.5
$ws.Range("L126").Value = 5580
$ws.Range("M126").Value = -2250.5
$ws.Range("N126").Value = -10520
